# Update cryptocurrency price/volume data per upstream feed refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.002.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.241.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.53%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.89"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.15%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.578.36"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.87"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.249.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.991.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0982"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.66%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.71"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +16.23%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.17"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.57"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.26"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0841"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.20%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0301"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.64"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.95"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.46"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.56"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.84%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.101"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.28%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.28"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -12.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.29"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.03%  "
